$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string header text (Volume/Number and Report Covering date range)
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# Crime-complaints table updates (rows 14-33)
$ws.Range("N14").Value = -44.444444444444
$ws.Range("C15").Value = "'0"
$ws.Range("L15").Value = 17.647058823529
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 71.428571428571
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 25.316455696202
$ws.Range("L16").Value = -2.941176470588
$ws.Range("M16").Value = -53.738317757009
$ws.Range("N16").Value = -84.482758620689
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 115.384615384615
$ws.Range("I17").Value = 219
$ws.Range("J17").Value = 196
$ws.Range("K17").Value = 11.734693877551
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 48.979591836734
$ws.Range("N17").Value = -13.779527559055
$ws.Range("C18").Value = "'0"
$ws.Range("F18").Value = 1
$ws.Range("H18").Value = -50
$ws.Range("L18").Value = -15.686274509803
$ws.Range("M18").Value = -82.157676348547
$ws.Range("N18").Value = -92.547660311958
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 265
$ws.Range("J19").Value = 263
$ws.Range("K19").Value = 0.760456273764
$ws.Range("L19").Value = 33.165829145728
$ws.Range("M19").Value = 17.256637168141
$ws.Range("N19").Value = -29.521276595744
$ws.Range("C20").Value = "'0"
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -16.666666666666
$ws.Range("J20").Value = 130
$ws.Range("K20").Value = -2.307692307692
$ws.Range("L20").Value = -14.76510067114
$ws.Range("M20").Value = 17.592592592592
$ws.Range("N20").Value = -93.662674650698
$ws.Range("C21").Value = 17
$ws.Range("E21").Value = 30.76923076923
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 25.862068965517
$ws.Range("I21").Value = 778
$ws.Range("J21").Value = 726
$ws.Range("K21").Value = 7.162534435261
$ws.Range("L21").Value = 4.993252361673
$ws.Range("M21").Value = -18.873826903024
$ws.Range("N21").Value = -79.948453608247
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 95
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = -5
$ws.Range("L23").Value = -17.391304347826
$ws.Range("M23").Value = 69.642857142857
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 400
$ws.Range("G24").Value = 34
$ws.Range("H24").Value = 44.117647058823
$ws.Range("I24").Value = 529
$ws.Range("J24").Value = 475
$ws.Range("K24").Value = 11.368421052631
$ws.Range("L24").Value = 6.653225806451
$ws.Range("M24").Value = 6.438631790744
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 144.444444444444
$ws.Range("I25").Value = 154
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = 38.738738738738
$ws.Range("L25").Value = 16.666666666666
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -18.518518518518
$ws.Range("I26").Value = 297
$ws.Range("J26").Value = 343
$ws.Range("K26").Value = -13.411078717201
$ws.Range("L26").Value = -12.389380530973
$ws.Range("M26").Value = -29.952830188679
$ws.Range("C27").Value = "'0"
$ws.Range("L27").Value = 10
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 30
$ws.Range("K28").Value = 3.448275862068
$ws.Range("L28").Value = 15.384615384615
$ws.Range("N29").Value = -55.555555555555
$ws.Range("N30").Value = -47.619047619047
$ws.Range("G33").Value = "'0"
$ws.Range("H33").Value = "'***.*"
